$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# "About" sheet
# ----------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Update the BNEF source block (rows 4-6) to reference the new article.
$about.Range("B4").Value = "Lithium-ion Battery Pack Prices Hit Record Low of `$139/kWh"
$about.Range("B5").Value = 2023
$about.Range("B6").Value = "https://about.bnef.com/blog/lithium-ion-battery-pack-prices-hit-record-low-of-139-kwh/#:~:text=Given%20this%2C%20BNEF%20expects%20average,and%20%2480%2FkWh%20in%202030."

# Remove the now-duplicated second source block (old rows 8-11).
$about.Rows("8:11").Delete()
# Remove the blank spacer row that used to separate the two blocks.
$about.Rows("7:7").Delete()
# Remove the old second numeric-conversion row (previously row 17).
$about.Rows("12:12").Delete()

# Replace the hyperlink on B6 with the updated article link.
$about.Hyperlinks.Add($about.Range("B6"), "https://about.bnef.com/blog/lithium-ion-battery-pack-prices-hit-record-low-of-139-kwh/", ":~:text=Given%20this%2C%20BNEF%20expects%20average,and%20%2480%2FkWh%20in%202030.") | Out-Null

# Replace the conversion-factor row (now row 11) with the new value.
$about.Range("A11").Value = "2023 to 2012"
$about.Range("B11").Value = 0.75350342301658668
$about.Range("B11").HorizontalAlignment = -4131

# ----------------------------------------------------------------------
# "BPP" sheet
# ----------------------------------------------------------------------
$bpp = $wb.Worksheets.Item("BPP")
$bpp.Range("B2").Formula = "=150*About!`$B`$11"
$bpp.Range("C2").Formula = "=161*About!`$B`$11"
$bpp.Range("D2").Formula = "=139*About!`$B`$11"
$bpp.Range("B3").Select() | Out-Null

# ----------------------------------------------------------------------
# "SYBPP" sheet
# ----------------------------------------------------------------------
$sybpp = $wb.Worksheets.Item("SYBPP")
$sybpp.Range("B2").Formula = "=160*About!B11"
$sybpp.Range("D13").Select() | Out-Null

# ----------------------------------------------------------------------
# Restore the "About" sheet as the active tab/selection, matching the
# workbook's last-saved UI state.
# ----------------------------------------------------------------------
$about.Select() | Out-Null
$about.Range("B21").Select() | Out-Null
